# Insert a new weekly price record as row 62 (Hortaliza / Agrícola del Norte
# S.A. de Arica - Albahaca), pushing the existing rows 62-74 down to 63-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 62..74 down to 63..75, leaving a blank row 62 to fill in.
$ws.Rows.Item(62).Insert()

# Populate the new row 62 with the new weekly record.
$ws.Cells.Item(62, 1).Value = 1
$ws.Cells.Item(62, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(62, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(62, 4).Value = 45173
$ws.Cells.Item(62, 5).Value = 15
$ws.Cells.Item(62, 6).Value = 100112052
$ws.Cells.Item(62, 7).Value = "Albahaca"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 420
$ws.Cells.Item(62, 11).Value = 900
$ws.Cells.Item(62, 12).Value = 1000
$ws.Cells.Item(62, 13).Value = 948
$ws.Cells.Item(62, 14).Value = "$/paquete"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 948
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"
